$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between B3 and C3
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 15

# Update C12 value
$ws.Range("C12").Value = 5

# Move selection/active cell to C12
$ws.Range("C12").Select()
